$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.980.13'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.465.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +13.13%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.466.10'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.59%  '
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.446'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.062.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.134'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000192'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("E16").Value = '  +3.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.018.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.468.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("E19").Value = '  +2.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("E23").Value = '  +3.95%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000120'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.10%  '
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +10.18%  '
$ws.Range("E31").Value = '  +1.27%  '
$ws.Range("E32").Value = '  +1.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.27%  '
$ws.Range("E35").Value = '  +12.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("E37").Value = '  +6.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0776'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.972.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.22%  '
$ws.Range("E42").Value = '  +5.53%  '
$ws.Range("E43").Value = '  +1.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.778'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.76'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.10%  '
$ws.Range("E47").Value = '  +3.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '314.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.51%  '
$ws.Range("E49").Value = '  +7.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.79%  '
$ws.Range("E51").Value = '  +4.56%  '
